$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $value)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "61.504.21"
Set-TextValue $ws "E2" "  -4.37%  "
Set-TextValue $ws "D3" "2.970.68"
Set-TextValue $ws "E3" "  -5.42%  "
Set-TextValue $ws "E4" "  +0.02%  "
Set-TextValue $ws "D5" "539.26"
Set-TextValue $ws "E5" "  -5.67%  "
Set-TextValue $ws "D6" "150.23"
Set-TextValue $ws "E6" "  -7.79%  "
Set-TextValue $ws "D7" "0.999"
Set-TextValue $ws "E7" "  +0.06%  "
Set-TextValue $ws "D8" "0.568"
Set-TextValue $ws "E8" "  -0.85%  "
Set-TextValue $ws "D9" "2.979.82"
Set-TextValue $ws "E9" "  -5.49%  "
Set-TextValue $ws "E10" "  -3.49%  "
Set-TextValue $ws "E11" "  -7.10%  "
Set-TextValue $ws "D12" "0.368"
Set-TextValue $ws "E12" "  -4.41%  "
Set-TextValue $ws "D13" "3.495.65"
Set-TextValue $ws "E13" "  -5.38%  "
Set-TextValue $ws "E14" "  -2.30%  "
Set-TextValue $ws "D15" "61.596.65"
Set-TextValue $ws "E15" "  -4.22%  "
Set-TextValue $ws "D16" "23.62"
Set-TextValue $ws "D17" "2.969.04"
Set-TextValue $ws "E17" "  -5.73%  "
Set-TextValue $ws "E18" "  -4.94%  "
Set-TextValue $ws "D19" "5.16"
Set-TextValue $ws "E19" "  -1.49%  "
Set-TextValue $ws "D20" "12.03"
Set-TextValue $ws "E20" "  -3.90%  "
Set-TextValue $ws "D21" "380.37"
Set-TextValue $ws "E21" "  -5.16%  "
Set-TextValue $ws "D22" "6.67"
Set-TextValue $ws "E22" "  -5.91%  "
Set-TextValue $ws "E23" "  +0.08%  "
Set-TextValue $ws "E24" "  -3.51%  "
Set-TextValue $ws "D25" "65.49"
Set-TextValue $ws "E25" "  -4.34%  "
Set-TextValue $ws "D26" "0.470"
Set-TextValue $ws "E26" "  -2.80%  "
Set-TextValue $ws "D27" "3.092.70"
Set-TextValue $ws "E27" "  -5.49%  "
Set-TextValue $ws "E28" "  -4.07%  "
Set-TextValue $ws "E29" "  +0.14%  "
Set-TextValue $ws "D30" "0.0₃0940"
Set-TextValue $ws "E30" "  -6.87%  "
Set-TextValue $ws "D31" "8.21"
Set-TextValue $ws "E31" "  -6.63%  "
Set-TextValue $ws "E32" "  +0.02%  "
Set-TextValue $ws "E33" "  -4.93%  "
Set-TextValue $ws "D34" "20.44"
Set-TextValue $ws "E34" "  -3.51%  "
Set-TextValue $ws "D35" "160.76"
Set-TextValue $ws "E35" "  -0.25%  "
Set-TextValue $ws "D36" "4.65"
Set-TextValue $ws "E36" "  -3.77%  "
Set-TextValue $ws "D37" "5.92"
Set-TextValue $ws "E37" "  -5.47%  "
Set-TextValue $ws "E38" "  -3.30%  "
Set-TextValue $ws "E39" "  -5.36%  "
Set-TextValue $ws "D40" "1.55"
Set-TextValue $ws "E40" "  -7.76%  "
Set-TextValue $ws "B41" "Filecoin"
Set-TextValue $ws "C41" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D41" "3.92"
Set-TextValue $ws "E41" "  -3.49%  "
Set-TextValue $ws "B42" "OKB"
Set-TextValue $ws "C42" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D42" "37.55"
Set-TextValue $ws "E42" "  -2.12%  "
Set-TextValue $ws "D43" "2.410.56"
Set-TextValue $ws "E43" "  -8.87%  "
Set-TextValue $ws "D44" "22.17"
Set-TextValue $ws "E44" "  -6.25%  "
Set-TextValue $ws "D45" "0.669"
Set-TextValue $ws "E45" "  -3.06%  "
Set-TextValue $ws "D46" "0.0590"
Set-TextValue $ws "E46" "  -3.55%  "
Set-TextValue $ws "E47" "  -6.73%  "
Set-TextValue $ws "E48" "  +0.11%  "
Set-TextValue $ws "D49" "0.0246"
Set-TextValue $ws "E49" "  -3.31%  "
Set-TextValue $ws "E50" "  -2.29%  "
Set-TextValue $ws "D51" "19.71"
Set-TextValue $ws "E51" "  -6.25%  "

Write-Host "Applied all changes"